# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the FFXIV Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 30050
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30050
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30050
$ws.Range("M44").ClearContents()  # was -4538
$ws.Range("N44").Value = -30974
$ws.Range("H45").Value = 51749.5
$ws.Range("I45").Value = 99999
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 299997
$ws.Range("L45").Value = 10500
$ws.Range("M45").Value = -299805
$ws.Range("N45").Value = -10884
$ws.Range("H137").Value = 1618.3208
$ws.Range("I137").Value = 1303.9524
$ws.Range("J137").Value = 2818.6365
$ws.Range("K137").Value = 3911.857199999999
$ws.Range("L137").Value = 8455.9095
$ws.Range("M137").Value = -1361.857199999999
$ws.Range("N137").Value = -13555.9095
$ws.Range("H138").Value = 2660.1729
$ws.Range("I138").Value = 1081.9348
$ws.Range("J138").Value = 4734.4287
$ws.Range("K138").Value = 3245.8044
$ws.Range("L138").Value = 14203.2861
$ws.Range("M138").Value = 1894.1956
$ws.Range("N138").Value = -24483.2861

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22264.95
$ws.Range("I32").Value = 25097.447
$ws.Range("K32").Value = 25097.447
$ws.Range("M32").Value = -24810.447
$ws.Range("H74").Value = 4948.4375
$ws.Range("I74").Value = 1892.12
$ws.Range("J74").Value = 15863.857
$ws.Range("K74").Value = 1892.12
$ws.Range("L74").Value = 15863.857
$ws.Range("M74").Value = -1018.12
$ws.Range("N74").Value = -17611.857
$ws.Range("H76").Value = 40288
$ws.Range("J76").Value = 40288
$ws.Range("L76").Value = 40288
$ws.Range("N76").Value = -40964
$ws.Range("H77").Value = 4948.4375
$ws.Range("I77").Value = 1892.12
$ws.Range("J77").Value = 15863.857
$ws.Range("K77").Value = 9460.599999999999
$ws.Range("L77").Value = 79319.285
$ws.Range("M77").Value = -5092.599999999999
$ws.Range("N77").Value = -88055.285
$ws.Range("H79").Value = 40288
$ws.Range("J79").Value = 40288
$ws.Range("L79").Value = 40288
$ws.Range("N79").Value = -42628
$ws.Range("H130").Value = 69000
$ws.Range("J130").Value = 69000
$ws.Range("L130").Value = 69000
$ws.Range("N130").Value = -79040

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 949.5
$ws.Range("I99").Value = 949.5
$ws.Range("K99").Value = 949.5
$ws.Range("M99").Value = 548.5
$ws.Range("H134").Value = 2564.5144
$ws.Range("I134").Value = 2291.6553
$ws.Range("K134").Value = 6874.965899999999
$ws.Range("M134").Value = -4339.965899999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1420
$ws.Range("I16").Value = 990.25
$ws.Range("J16").Value = 1527.4375
$ws.Range("K16").Value = 990.25
$ws.Range("L16").Value = 1527.4375
$ws.Range("M16").Value = -703.25
$ws.Range("N16").Value = -2101.4375
$ws.Range("H31").Value = 2364.5
$ws.Range("I31").Value = 1811.6897
$ws.Range("J31").Value = 3208.2632
$ws.Range("K31").Value = 1811.6897
$ws.Range("L31").Value = 3208.2632
$ws.Range("M31").Value = -1516.6897
$ws.Range("N31").Value = -3798.2632
$ws.Range("H34").Value = 2364.5
$ws.Range("I34").Value = 1811.6897
$ws.Range("J34").Value = 3208.2632
$ws.Range("K34").Value = 1811.6897
$ws.Range("L34").Value = 3208.2632
$ws.Range("M34").Value = -1609.6897
$ws.Range("N34").Value = -3612.2632
$ws.Range("H58").Value = 2844249.5
$ws.Range("I58").Value = 4547968.5
$ws.Range("J58").Value = 4717.8335
$ws.Range("K58").Value = 4547968.5
$ws.Range("L58").Value = 4717.8335
$ws.Range("M58").Value = -4547765.5
$ws.Range("N58").Value = -5123.8335
$ws.Range("H113").Value = 1420
$ws.Range("I113").Value = 990.25
$ws.Range("J113").Value = 1527.4375
$ws.Range("K113").Value = 990.25
$ws.Range("L113").Value = 1527.4375
$ws.Range("M113").Value = 1179.75
$ws.Range("N113").Value = -5867.4375
$ws.Range("H114").Value = 63684
$ws.Range("J114").Value = 63684
$ws.Range("L114").Value = 63684
$ws.Range("N114").Value = -72362
$ws.Range("H136").Value = 2844249.5
$ws.Range("I136").Value = 4547968.5
$ws.Range("J136").Value = 4717.8335
$ws.Range("K136").Value = 13643905.5
$ws.Range("L136").Value = 14153.5005
$ws.Range("M136").Value = -13641355.5
$ws.Range("N136").Value = -19253.5005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 605
$ws.Range("I44").Value = 605
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1815
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -1417
$ws.Range("N44").ClearContents()  # was -3286.9999
$ws.Range("H120").Value = 8632.5
$ws.Range("I120").Value = 8676.666999999999
$ws.Range("K120").Value = 26030.001
$ws.Range("M120").Value = -21192.001
$ws.Range("H131").Value = 15010.853
$ws.Range("I131").Value = 1153.3334
$ws.Range("J131").Value = 16351.903
$ws.Range("K131").Value = 3460.0002
$ws.Range("L131").Value = 49055.709
$ws.Range("M131").Value = 1579.9998
$ws.Range("N131").Value = -59135.709
$ws.Range("H137").Value = 17762.361
$ws.Range("I137").Value = 1536.3572
$ws.Range("J137").Value = 28088
$ws.Range("K137").Value = 4609.071599999999
$ws.Range("L137").Value = 84264
$ws.Range("M137").Value = 490.9284000000007
$ws.Range("N137").Value = -94464

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 48000
$ws.Range("J74").Value = 48000
$ws.Range("L74").Value = 48000
$ws.Range("N74").Value = -49872
$ws.Range("H77").Value = 48000
$ws.Range("J77").Value = 48000
$ws.Range("L77").Value = 144000
$ws.Range("N77").Value = -153360
$ws.Range("H107").Value = 709.1111
$ws.Range("I107").Value = 286
$ws.Range("J107").Value = 1238
$ws.Range("K107").Value = 286
$ws.Range("L107").Value = 1238
$ws.Range("M107").Value = 1634
$ws.Range("N107").Value = -5078

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3628.9517
$ws.Range("I132").Value = 3634.8293
$ws.Range("J132").Value = 3617.476
$ws.Range("K132").Value = 10904.4879
$ws.Range("L132").Value = 10852.428
$ws.Range("M132").Value = -8374.4879
$ws.Range("N132").Value = -15912.428

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2639.4443
$ws.Range("I132").Value = 1265.5834
$ws.Range("J132").Value = 3738.5334
$ws.Range("K132").Value = 3796.7502
$ws.Range("L132").Value = 11215.6002
$ws.Range("M132").Value = -1266.7502
$ws.Range("N132").Value = -16275.6002
$ws.Range("H136").Value = 5170.9395
$ws.Range("I136").Value = 4449.794
$ws.Range("J136").Value = 5937.1562
$ws.Range("K136").Value = 13349.382
$ws.Range("L136").Value = 17811.4686
$ws.Range("M136").Value = -10799.382
$ws.Range("N136").Value = -22911.4686

